$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2333.3333
$ws.Range("J32").Value = 2500
$ws.Range("L32").Value = 2500
$ws.Range("N32").Value = -3152
$ws.Range("H98").Value = 3945.5417
$ws.Range("I98").Value = 4627.1665
$ws.Range("K98").Value = 4627.1665
$ws.Range("M98").Value = -3129.1665
$ws.Range("H106").Value = 9752.9375
$ws.Range("I106").Value = 10002.733
$ws.Range("K106").Value = 10002.733
$ws.Range("M106").Value = -9371.733
$ws.Range("H118").Value = 729.6667
$ws.Range("I118").Value = 729.6667
$ws.Range("K118").Value = 2189.0001
$ws.Range("M118").Value = -532.0001000000002
$ws.Range("H122").Value = 3945.5417
$ws.Range("I122").Value = 4627.1665
$ws.Range("K122").Value = 13881.4995
$ws.Range("M122").Value = -11431.4995
$ws.Range("H137").Value = 1070.6173
$ws.Range("I137").Value = 894.3721
$ws.Range("J137").Value = 1270.0526
$ws.Range("K137").Value = 2683.1163
$ws.Range("L137").Value = 3810.1578
$ws.Range("M137").Value = -133.1163000000001
$ws.Range("N137").Value = -8910.157800000001
$ws.Range("H138").Value = 596655.25
$ws.Range("I138").Value = 711.0213
$ws.Range("J138").Value = 1673939
$ws.Range("K138").Value = 2133.0639
$ws.Range("L138").Value = 5021817
$ws.Range("M138").Value = 3006.9361
$ws.Range("N138").Value = -5032097
$ws.Range("H139").Value = 34520
$ws.Range("J139").Value = 34520
$ws.Range("L139").Value = 34520
$ws.Range("N139").Value = -44800

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4976.7827
$ws.Range("I32").Value = 4691.636
$ws.Range("J32").Value = 11250
$ws.Range("K32").Value = 4691.636
$ws.Range("L32").Value = 11250
$ws.Range("M32").Value = -4404.636
$ws.Range("N32").Value = -11824
$ws.Range("H110").Value = 1958.8235
$ws.Range("I110").Value = 1506.3
$ws.Range("J110").Value = 2605.2856
$ws.Range("K110").Value = 1506.3
$ws.Range("L110").Value = 2605.2856
$ws.Range("M110").Value = 538.7
$ws.Range("N110").Value = -6695.2856
$ws.Range("H132").Value = 2281.4187
$ws.Range("I132").Value = 2249.6333
$ws.Range("J132").Value = 2354.7693
$ws.Range("K132").Value = 6748.8999
$ws.Range("L132").Value = 7064.3079
$ws.Range("M132").Value = -4218.8999
$ws.Range("N132").Value = -12124.3079
$ws.Range("H139").Value = 31024.9
$ws.Range("J139").Value = 31024.9
$ws.Range("L139").Value = 31024.9
$ws.Range("N139").Value = -41304.9

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2579.8
$ws.Range("I86").Value = 3014.5186
$ws.Range("J86").Value = 1676.9231
$ws.Range("K86").Value = 3014.5186
$ws.Range("L86").Value = 1676.9231
$ws.Range("M86").Value = -1891.5186
$ws.Range("N86").Value = -3922.9231
$ws.Range("H89").Value = 2579.8
$ws.Range("I89").Value = 3014.5186
$ws.Range("J89").Value = 1676.9231
$ws.Range("K89").Value = 15072.593
$ws.Range("L89").Value = 8384.6155
$ws.Range("M89").Value = -9456.592999999999
$ws.Range("N89").Value = -19616.6155
$ws.Range("H105").Value = 36069084
$ws.Range("I105").Value = 40397160
$ws.Range("J105").Value = 1803.3334
$ws.Range("K105").Value = 40397160
$ws.Range("L105").Value = 1803.3334
$ws.Range("M105").Value = -40395413
$ws.Range("N105").Value = -5297.3334
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 30000
$ws.Range("N120").Value = -39676
$ws.Range("H134").Value = 4485.7075
$ws.Range("I134").Value = 1367.1212
$ws.Range("J134").Value = 17349.875
$ws.Range("K134").Value = 4101.363600000001
$ws.Range("L134").Value = 52049.625
$ws.Range("M134").Value = -1566.363600000001
$ws.Range("N134").Value = -57119.625

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6252316
$ws.Range("I62").Value = 2378.75
$ws.Range("J62").Value = 50001876
$ws.Range("K62").Value = 2378.75
$ws.Range("L62").Value = 50001876
$ws.Range("M62").Value = -1754.75
$ws.Range("N62").Value = -50003124
$ws.Range("H65").Value = 6252316
$ws.Range("I65").Value = 2378.75
$ws.Range("J65").Value = 50001876
$ws.Range("K65").Value = 11893.75
$ws.Range("L65").Value = 250009380
$ws.Range("M65").Value = -8773.75
$ws.Range("N65").Value = -250015620
$ws.Range("H112").Value = 37175.5
$ws.Range("J112").Value = 37175.5
$ws.Range("L112").Value = 37175.5
$ws.Range("N112").Value = -40129.5
$ws.Range("H120").Value = 20000
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 20000
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 20000
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -27258
$ws.Range("H132").Value = 3326.6785
$ws.Range("I132").Value = 3217.4583
$ws.Range("J132").Value = 3982
$ws.Range("K132").Value = 9652.374899999999
$ws.Range("L132").Value = 11946
$ws.Range("M132").Value = -7122.374899999999
$ws.Range("N132").Value = -17006

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 183333.33
$ws.Range("J37").Value = 183333.33
$ws.Range("L37").Value = 549999.99
$ws.Range("N37").Value = -550223.99

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3156.6667
$ws.Range("J80").Value = 4242
$ws.Range("L80").Value = 4242
$ws.Range("N80").Value = -6238
$ws.Range("H83").Value = 3156.6667
$ws.Range("J83").Value = 4242
$ws.Range("L83").Value = 21210
$ws.Range("N83").Value = -31194
$ws.Range("H102").Value = 1488.1818
$ws.Range("I102").Value = 1301.4286
$ws.Range("J102").Value = 1815
$ws.Range("K102").Value = 1301.4286
$ws.Range("L102").Value = 1815
$ws.Range("M102").Value = 320.5714
$ws.Range("N102").Value = -5059
$ws.Range("H107").Value = 763.8570999999999
$ws.Range("J107").Value = 449
$ws.Range("L107").Value = 449
$ws.Range("N107").Value = -4289
$ws.Range("H126").Value = 1848.9615
$ws.Range("I126").Value = 1593
$ws.Range("J126").Value = 2258.5
$ws.Range("K126").Value = 4779
$ws.Range("L126").Value = 6775.5
$ws.Range("M126").Value = -2309
$ws.Range("N126").Value = -11715.5
$ws.Range("H132").Value = 1922.9117
$ws.Range("I132").Value = 1662.72
$ws.Range("K132").Value = 4988.16
$ws.Range("M132").Value = -2458.16
$ws.Range("H135").Value = 37651.8
$ws.Range("J135").Value = 34564.75
$ws.Range("L135").Value = 34564.75
$ws.Range("N135").Value = -44704.75
$ws.Range("H140").Value = 33840
$ws.Range("J140").Value = 33840
$ws.Range("L140").Value = 33840
$ws.Range("N140").Value = -44200

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19128.05
$ws.Range("I132").Value = 1465.1177
$ws.Range("J132").Value = 44150.543
$ws.Range("K132").Value = 4395.3531
$ws.Range("L132").Value = 132451.629
$ws.Range("M132").Value = -1865.3531
$ws.Range("N132").Value = -137511.629
$ws.Range("H135").Value = 37015.875
$ws.Range("J135").Value = 37015.875
$ws.Range("L135").Value = 37015.875
$ws.Range("N135").Value = -47155.875
$ws.Range("H141").Value = 70215
$ws.Range("J141").Value = 70215
$ws.Range("L141").Value = 70215
$ws.Range("N141").Value = -80575

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 70004
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 70004
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 70004
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -70284
$ws.Range("H126").Value = 40001416
$ws.Range("I126").Value = 52632816
$ws.Range("J126").Value = 1983.3334
$ws.Range("K126").Value = 157898448
$ws.Range("L126").Value = 5950.0002
$ws.Range("M126").Value = -157895978
$ws.Range("N126").Value = -10890.0002
$ws.Range("H132").Value = 1722.9672
$ws.Range("I132").Value = 1661.8909
$ws.Range("K132").Value = 4985.6727
$ws.Range("M132").Value = -2455.6727
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
